# Create test case for login error
# Renames Sheet1 -> valid_Account, adds a new sheet invalid_Account with
# negative/edge-case login data, and updates selections/active tab to match.

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename to valid_Account, clear the old "tab selected" /
#     single-cell selection so selection becomes the whole used range ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "valid_Account"
$ws1.Range("A1:B3").Select() | Out-Null

# --- Add the new sheet right after valid_Account ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "invalid_Account"

# Header row (reuses the same shared strings as valid_Account)
$ws2.Range("A1").Value = "userName"
$ws2.Range("B1").Value = "password"

# Row 2 - blank/whitespace username with a valid-looking password
$ws2.Range("B2").Value = "Test@1"
$ws2.Range("A2").Value = "    "

# Row 4 - whitespace-only username & password
$ws2.Range("A4").Value = "      "
$ws2.Range("B4").Value = "       "

# Row 5 - special characters
$ws2.Range("A5").Value = "))@#*@#@_"
$ws2.Range("B5").Value = ". "

# Row 3 - new username reusing the existing valid password / hyperlink
$ws2.Range("A3").Value = "binhvh3"
$ws2.Range("B3").Value = "aA@123456"

# Hyperlinks: B3 (existing valid password) gets relationship id 1,
# then B2 (new password) gets relationship id 2
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:aA@123456") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:Test@1") | Out-Null
$ws2.Range("B3").Style = "Hyperlink"
$ws2.Range("B2").Style = "Hyperlink"

# Row 6 - SQL-injection-style username, numeric password
$ws2.Range("A6").Value = "105 OR 1=1"
$ws2.Range("B6").Value = 1

# Row 7 - SQL-injection-style username/password pair (identical)
$ws2.Range("A7").Value = " or ""="
$ws2.Range("B7").Value = " or ""="

# Column sizing similar to the source sheet
$ws2.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# This is the sheet/tab the workbook opens on
$ws2.Range("X7").Select() | Out-Null

Write-Output "Workbook now has sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
